# Updated TPM-derived ligand/receptor expression + specificity + edge-weight
# values for the Cd274-Cd80 LR-pair sheet (NATMI lrc2p output), per commit
# "update scripts wuth new tpm". Columns G:J (ligand expr/specificity),
# M:P (receptor expr/specificity) and Q:T (edge weight/specificity) are
# recomputed from the new TPM values; A:F and K:L are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 10.43877666666667
$ws.Range("H2").Value = 31.31633
$ws.Range("I2").Value = 0.478485024444405
$ws.Range("J2").Value = 0.4784850244444049
$ws.Range("M2").Value = 1.152905666666667
$ws.Range("N2").Value = 3.458717
$ws.Range("O2").Value = 0.06522949989114324
$ws.Range("P2").Value = 0.06522949989114325
$ws.Range("Q2").Value = 12.03492477206778
$ws.Range("R2").Value = 108.31432294861
$ws.Range("S2").Value = 0.03121133884990998
$ws.Range("T2").Value = 0.03121133884990998

# Row 3
$ws.Range("G3").Value = 10.43877666666667
$ws.Range("H3").Value = 31.31633
$ws.Range("I3").Value = 0.478485024444405
$ws.Range("J3").Value = 0.4784850244444049
$ws.Range("O3").Value = 0.2252875952949142
$ws.Range("P3").Value = 0.2252875952949143
$ws.Range("Q3").Value = 41.56584468651556
$ws.Range("R3").Value = 374.09260217864
$ws.Range("S3").Value = 0.1077967405417083
$ws.Range("T3").Value = 0.1077967405417082

# Row 4
$ws.Range("G4").Value = 10.43877666666667
$ws.Range("H4").Value = 31.31633
$ws.Range("I4").Value = 0.478485024444405
$ws.Range("J4").Value = 0.4784850244444049
$ws.Range("M4").Value = 1.517768666666667
$ws.Range("N4").Value = 4.553306
$ws.Range("O4").Value = 0.08587284626968379
$ws.Range("P4").Value = 0.08587284626968379
$ws.Range("Q4").Value = 15.84364814299778
$ws.Range("R4").Value = 142.59283328698
$ws.Range("S4").Value = 0.04108887094646028
$ws.Range("T4").Value = 0.04108887094646027

# Row 5
$ws.Range("G5").Value = 10.43877666666667
$ws.Range("H5").Value = 31.31633
$ws.Range("I5").Value = 0.478485024444405
$ws.Range("J5").Value = 0.4784850244444049
$ws.Range("M5").Value = 11.02206166666667
$ws.Range("N5").Value = 33.066185
$ws.Range("O5").Value = 0.6236100585442587
$ws.Range("P5").Value = 0.6236100585442588
$ws.Range("Q5").Value = 115.0568401445611
$ws.Range("R5").Value = 1035.51156130105
$ws.Range("S5").Value = 0.2983880741063265
$ws.Range("T5").Value = 0.2983880741063264

# Row 6
$ws.Range("I6").Value = 0.07482651662844755
$ws.Range("J6").Value = 0.07482651662844754
$ws.Range("M6").Value = 1.152905666666667
$ws.Range("N6").Value = 3.458717
$ws.Range("O6").Value = 0.06522949989114324
$ws.Range("P6").Value = 0.06522949989114325
$ws.Range("Q6").Value = 1.882047404983889
$ws.Range("R6").Value = 16.938426644855
$ws.Range("S6").Value = 0.004880896258269947
$ws.Range("T6").Value = 0.004880896258269948

# Row 7
$ws.Range("I7").Value = 0.07482651662844755
$ws.Range("J7").Value = 0.07482651662844754
$ws.Range("O7").Value = 0.2252875952949142
$ws.Range("P7").Value = 0.2252875952949143
$ws.Range("S7").Value = 0.01685748599551786
$ws.Range("T7").Value = 0.01685748599551786

# Row 8
$ws.Range("I8").Value = 0.07482651662844755
$ws.Range("J8").Value = 0.07482651662844754
$ws.Range("M8").Value = 1.517768666666667
$ws.Range("N8").Value = 4.553306
$ws.Range("O8").Value = 0.08587284626968379
$ws.Range("P8").Value = 0.08587284626968379
$ws.Range("Q8").Value = 2.477663752598889
$ws.Range("R8").Value = 22.29897377339
$ws.Range("S8").Value = 0.006425565959330614
$ws.Range("T8").Value = 0.006425565959330613

# Row 9
$ws.Range("I9").Value = 0.07482651662844755
$ws.Range("J9").Value = 0.07482651662844754
$ws.Range("M9").Value = 11.02206166666667
$ws.Range("N9").Value = 33.066185
$ws.Range("O9").Value = 0.6236100585442587
$ws.Range("P9").Value = 0.6236100585442588
$ws.Range("Q9").Value = 17.99283597703056
$ws.Range("R9").Value = 161.935523793275
$ws.Range("S9").Value = 0.04666256841532912
$ws.Range("T9").Value = 0.04666256841532912

# Row 10
$ws.Range("G10").Value = 0.4778236666666666
$ws.Range("H10").Value = 1.433471
$ws.Range("I10").Value = 0.02190213241702797
$ws.Range("J10").Value = 0.02190213241702797
$ws.Range("M10").Value = 1.152905666666667
$ws.Range("N10").Value = 3.458717
$ws.Range("O10").Value = 0.06522949989114324
$ws.Range("P10").Value = 0.06522949989114325
$ws.Range("Q10").Value = 0.5508856129674444
$ws.Range("R10").Value = 4.957970516707
$ws.Range("S10").Value = 0.001428665144112331
$ws.Range("T10").Value = 0.001428665144112331

# Row 11
$ws.Range("G11").Value = 0.4778236666666666
$ws.Range("H11").Value = 1.433471
$ws.Range("I11").Value = 0.02190213241702797
$ws.Range("J11").Value = 0.02190213241702797
$ws.Range("O11").Value = 0.2252875952949142
$ws.Range("P11").Value = 0.2252875952949143
$ws.Range("Q11").Value = 1.902631405040889
$ws.Range("R11").Value = 17.123682645368
$ws.Range("S11").Value = 0.00493427874406302
$ws.Range("T11").Value = 0.00493427874406302

# Row 12
$ws.Range("G12").Value = 0.4778236666666666
$ws.Range("H12").Value = 1.433471
$ws.Range("I12").Value = 0.02190213241702797
$ws.Range("J12").Value = 0.02190213241702797
$ws.Range("M12").Value = 1.517768666666667
$ws.Range("N12").Value = 4.553306
$ws.Range("O12").Value = 0.08587284626968379
$ws.Range("P12").Value = 0.08587284626968379
$ws.Range("Q12").Value = 0.7252257894584445
$ws.Range("R12").Value = 6.527032105126
$ws.Range("S12").Value = 0.001880798450025701
$ws.Range("T12").Value = 0.001880798450025701

# Row 13
$ws.Range("G13").Value = 0.4778236666666666
$ws.Range("H13").Value = 1.433471
$ws.Range("I13").Value = 0.02190213241702797
$ws.Range("J13").Value = 0.02190213241702797
$ws.Range("M13").Value = 11.02206166666667
$ws.Range("N13").Value = 33.066185
$ws.Range("O13").Value = 0.6236100585442587
$ws.Range("P13").Value = 0.6236100585442588
$ws.Range("Q13").Value = 5.266601919792778
$ws.Range("R13").Value = 47.399417278135
$ws.Range("S13").Value = 0.01365839007882692
$ws.Range("T13").Value = 0.01365839007882692

# Row 14
$ws.Range("G14").Value = 9.267269333333333
$ws.Range("H14").Value = 27.801808
$ws.Range("I14").Value = 0.4247863265101195
$ws.Range("J14").Value = 0.4247863265101195
$ws.Range("M14").Value = 1.152905666666667
$ws.Range("N14").Value = 3.458717
$ws.Range("O14").Value = 0.06522949989114324
$ws.Range("P14").Value = 0.06522949989114325
$ws.Range("Q14").Value = 10.68428732892622
$ws.Range("R14").Value = 96.15858596033601
$ws.Range("S14").Value = 0.02770859963885098
$ws.Range("T14").Value = 0.02770859963885099

# Row 15
$ws.Range("G15").Value = 9.267269333333333
$ws.Range("H15").Value = 27.801808
$ws.Range("I15").Value = 0.4247863265101195
$ws.Range("J15").Value = 0.4247863265101195
$ws.Range("O15").Value = 0.2252875952949142
$ws.Range("P15").Value = 0.2252875952949143
$ws.Range("Q15").Value = 36.90105556214044
$ws.Range("R15").Value = 332.109500059264
$ws.Range("S15").Value = 0.09569909001362512
$ws.Range("T15").Value = 0.09569909001362513

# Row 16
$ws.Range("G16").Value = 9.267269333333333
$ws.Range("H16").Value = 27.801808
$ws.Range("I16").Value = 0.4247863265101195
$ws.Range("J16").Value = 0.4247863265101195
$ws.Range("M16").Value = 1.517768666666667
$ws.Range("N16").Value = 4.553306
$ws.Range("O16").Value = 0.08587284626968379
$ws.Range("P16").Value = 0.08587284626968379
$ws.Range("Q16").Value = 14.06557101969422
$ws.Range("R16").Value = 126.590139177248
$ws.Range("S16").Value = 0.0364776109138672
$ws.Range("T16").Value = 0.0364776109138672

# Row 17
$ws.Range("G17").Value = 9.267269333333333
$ws.Range("H17").Value = 27.801808
$ws.Range("I17").Value = 0.4247863265101195
$ws.Range("J17").Value = 0.4247863265101195
$ws.Range("M17").Value = 11.02206166666667
$ws.Range("N17").Value = 33.066185
$ws.Range("O17").Value = 0.6236100585442587
$ws.Range("P17").Value = 0.6236100585442588
$ws.Range("Q17").Value = 102.1444140736089
$ws.Range("R17").Value = 919.2997266624801
$ws.Range("S17").Value = 0.2649010259437762
$ws.Range("T17").Value = 0.2649010259437763
